$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (currently sitting in the
#    empty paragraph at the very end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Drop the trailing "s" from "runtimes." -> "runtime." in the
#    "most important adjustments" bullet.
$d.Content.Find.Execute("runtimes.", $true, $false, $false, $false, $false, $true, 1, $false, "runtime.", 2) | Out-Null

# 3. Re-insert the "_GoBack" bookmark right before the trailing period of
#    that same sentence, which splits the run in two (mirroring a real
#    Word edit session that left the cursor there).
$found = $d.Content
$found.Find.Execute("runtime.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkPos = $found.End - 1
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
